$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Find the last used row based on column A
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

for ($r = 1; $r -le $lastRow; $r++) {
    $cVal = $ws.Cells.Item($r, 3).Value2
    $dVal = $ws.Cells.Item($r, 4).Value2
    $ws.Cells.Item($r, 3).Value2 = $dVal
    $ws.Cells.Item($r, 4).Value2 = $cVal
}
